{"js": "// Compte rendu V 1.2\n// Fix a typo in the \"capteur / enfant\" paragraph: \"\u00e7a main\" -> \"sa main\".\n// (\"\u00e7a\" only occurs once in the whole document body, inside the sentence\n//  \"...si un enfant tente de saisir le produit nocif, \u00e7a main \u00e9tant plus\n//  petite qu'une main d'adulte...\" \u2014 it should read \"sa main\".)\n\nconst body = context.document.body;\n\n// Scope the search narrowly to the unique phrase so we don't accidentally\n// touch anything else, then load the match.\nconst results = body.search(\"\u00e7a main\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items,text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error('Expected to find \"\u00e7a main\" in the document body, found none.');\n}\n\n// Replace \"\u00e7a\" -> \"sa\" while keeping \" main\" intact, preserving the run's\n// existing character formatting (font, size, color, etc.).\nconst target = results.items[0];\ntarget.insertText(\"sa main\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Compte rendu V 1.2\n# Fix a typo in the \"capteur / enfant\" paragraph: \"\u00e7a main\" -> \"sa main\".\n# (\"\u00e7a\" only occurs once in the whole document, inside the sentence\n#  \"...si un enfant tente de saisir le produit nocif, \u00e7a main \u00e9tant plus\n#  petite qu'une main d'adulte...\" \u2014 it should read \"sa main\".)\n\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"\u00e7a main\"\n$find.Replacement.Text = \"sa main\"\n\n# wdFindContinue = 1, wdReplaceOne = 1 (Execute's Replace arg)\n$find.Execute([ref]$find.Text, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$find.Replacement.Text, [ref]1)\n"}
